$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 336, shifting existing rows 336:411 down to 337:412
$ws.Rows.Item(336).Insert()

# Populate the new row 336 with its data (same categorical columns as the
# row that used to occupy position 336, but new observation values)
$ws.Cells.Item(336, 1).Value = 3
$ws.Cells.Item(336, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(336, 3).Value = "Coquimbo"
$ws.Cells.Item(336, 4).Value = 44798
$ws.Cells.Item(336, 5).Value = 5
$ws.Cells.Item(336, 6).Value = 100112040
$ws.Cells.Item(336, 7).Value = "Cilantro"
$ws.Cells.Item(336, 8).Value = "Sin especificar"
$ws.Cells.Item(336, 9).Value = "Primera"
$ws.Cells.Item(336, 10).Value = 190
$ws.Cells.Item(336, 11).Value = 4000
$ws.Cells.Item(336, 12).Value = 4500
$ws.Cells.Item(336, 13).Value = 4289
$ws.Cells.Item(336, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(336, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(336, 16).Value = 1430
$ws.Cells.Item(336, 17).Value = 3
$ws.Cells.Item(336, 18).Value = "Hortaliza"
